# chore: update Sheets via scheduled runner
# Refresh cached market-price / profit figures (columns H-N: currentAveragePrice,
# currentAveragePriceNQ, currentAveragePriceHQ, LevePriceNQ, LevePriceHQ,
# LeveProfitNQ, LeveProfitHQ) for the affected leve rows across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 219.3077
$ws.Cells.Item(12, 9).Value = 154.5
$ws.Cells.Item(12, 11).Value = 154.5
$ws.Cells.Item(12, 13).Value = 15.5

$ws.Cells.Item(86, 8).Value = 5512.3335
$ws.Cells.Item(86, 9).Value = 5528.143
$ws.Cells.Item(86, 11).Value = 5528.143
$ws.Cells.Item(86, 13).Value = -4405.143

$ws.Cells.Item(89, 8).Value = 5512.3335
$ws.Cells.Item(89, 9).Value = 5528.143
$ws.Cells.Item(89, 11).Value = 27640.715
$ws.Cells.Item(89, 13).Value = -22024.715

$ws.Cells.Item(131, 8).Value = 90931640
$ws.Cells.Item(131, 9).Value = 166668830
$ws.Cells.Item(131, 11).Value = 500006490
$ws.Cells.Item(131, 13).Value = -500001450

$ws.Cells.Item(135, 8).Value = 5218.913
$ws.Cells.Item(135, 9).Value = 1153.9231
$ws.Cells.Item(135, 11).Value = 10385.3079
$ws.Cells.Item(135, 13).Value = -7850.3079

$ws.Cells.Item(137, 8).Value = 15877481
$ws.Cells.Item(137, 9).Value = 38463236
$ws.Cells.Item(137, 10).Value = 6409.243
$ws.Cells.Item(137, 11).Value = 115389708
$ws.Cells.Item(137, 12).Value = 19227.729
$ws.Cells.Item(137, 13).Value = -115387158
$ws.Cells.Item(137, 14).Value = -24327.729

$ws.Cells.Item(138, 8).Value = 2884.3333
$ws.Cells.Item(138, 10).Value = 3357.4482
$ws.Cells.Item(138, 12).Value = 10072.3446
$ws.Cells.Item(138, 14).Value = -20352.3446

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 144681.53
$ws.Cells.Item(32, 9).Value = 222763.44
$ws.Cells.Item(32, 10).Value = 20827.482
$ws.Cells.Item(32, 11).Value = 222763.44
$ws.Cells.Item(32, 12).Value = 20827.482
$ws.Cells.Item(32, 13).Value = -222476.44
$ws.Cells.Item(32, 14).Value = -21401.482

$ws.Cells.Item(61, 8).Value = 2782203.5
$ws.Cells.Item(61, 9).Value = 4454.8125
$ws.Cells.Item(61, 11).Value = 4454.8125
$ws.Cells.Item(61, 13).Value = -4242.8125

$ws.Cells.Item(74, 8).Value = 1188757.4
$ws.Cells.Item(74, 9).Value = 1427107.6
$ws.Cells.Item(74, 11).Value = 1427107.6
$ws.Cells.Item(74, 13).Value = -1426233.6

$ws.Cells.Item(77, 8).Value = 1188757.4
$ws.Cells.Item(77, 9).Value = 1427107.6
$ws.Cells.Item(77, 11).Value = 7135538
$ws.Cells.Item(77, 13).Value = -7131170

$ws.Cells.Item(97, 8).Value = 142861250
$ws.Cells.Item(97, 9).Value = 4125.75
$ws.Cells.Item(97, 11).Value = 4125.75
$ws.Cells.Item(97, 13).Value = -3629.75

$ws.Cells.Item(132, 8).Value = 761076.0600000001
$ws.Cells.Item(132, 9).Value = 809906.8
$ws.Cells.Item(132, 11).Value = 2429720.4
$ws.Cells.Item(132, 13).Value = -2427190.4

$ws.Cells.Item(133, 8).Value = 70248.5
$ws.Cells.Item(133, 10).Value = 70248.5
$ws.Cells.Item(133, 12).Value = 70248.5
$ws.Cells.Item(133, 14).Value = -75308.5

$ws.Cells.Item(134, 8).Value = 46540.2
$ws.Cells.Item(134, 10).Value = 46540.2
$ws.Cells.Item(134, 12).Value = 46540.2
$ws.Cells.Item(134, 14).Value = -56680.2

$ws.Cells.Item(136, 8).Value = 2782203.5
$ws.Cells.Item(136, 9).Value = 4454.8125
$ws.Cells.Item(136, 11).Value = 13364.4375
$ws.Cells.Item(136, 13).Value = -10814.4375

$ws.Cells.Item(138, 8).Value = 116332.664
$ws.Cells.Item(138, 10).Value = 116332.664
$ws.Cells.Item(138, 12).Value = 116332.664
$ws.Cells.Item(138, 14).Value = -126612.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 57616.2
$ws.Cells.Item(20, 9).Value = 66917.35000000001
$ws.Cells.Item(20, 11).Value = 66917.35000000001
$ws.Cells.Item(20, 13).Value = -66670.35000000001

$ws.Cells.Item(94, 8).Value = 323835.25
$ws.Cells.Item(94, 9).Value = 12312.214
$ws.Cells.Item(94, 11).Value = 12312.214
$ws.Cells.Item(94, 13).Value = -11861.214

$ws.Cells.Item(107, 8).Value = 1094.55
$ws.Cells.Item(107, 9).Value = 814.7692
$ws.Cells.Item(107, 10).Value = 1614.1428
$ws.Cells.Item(107, 11).Value = 814.7692
$ws.Cells.Item(107, 12).Value = 1614.1428
$ws.Cells.Item(107, 13).Value = 1105.2308
$ws.Cells.Item(107, 14).Value = -5454.1428

$ws.Cells.Item(134, 8).Value = 2420988.5
$ws.Cells.Item(134, 9).Value = 3721.709
$ws.Cells.Item(134, 11).Value = 11165.127
$ws.Cells.Item(134, 13).Value = -8630.127

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 54.565216
$ws.Cells.Item(7, 9).Value = 67.5625
$ws.Cells.Item(7, 10).Value = 24.857143
$ws.Cells.Item(7, 11).Value = 67.5625
$ws.Cells.Item(7, 12).Value = 24.857143
$ws.Cells.Item(7, 13).Value = 45.4375
$ws.Cells.Item(7, 14).Value = -250.857143

$ws.Cells.Item(58, 8).Value = 2612579.5
$ws.Cells.Item(58, 9).Value = 5216.3687
$ws.Cells.Item(58, 10).Value = 6423341
$ws.Cells.Item(58, 11).Value = 5216.3687
$ws.Cells.Item(58, 12).Value = 6423341
$ws.Cells.Item(58, 13).Value = -5013.3687
$ws.Cells.Item(58, 14).Value = -6423747

$ws.Cells.Item(62, 8).Value = 4573.8
$ws.Cells.Item(62, 10).Value = 4469
$ws.Cells.Item(62, 12).Value = 4469
$ws.Cells.Item(62, 14).Value = -5717

$ws.Cells.Item(65, 8).Value = 4573.8
$ws.Cells.Item(65, 10).Value = 4469
$ws.Cells.Item(65, 12).Value = 22345
$ws.Cells.Item(65, 14).Value = -28585

$ws.Cells.Item(107, 8).Value = 720
$ws.Cells.Item(107, 9).Value = 605.4
$ws.Cells.Item(107, 10).Value = 1006.5
$ws.Cells.Item(107, 11).Value = 605.4
$ws.Cells.Item(107, 12).Value = 1006.5
$ws.Cells.Item(107, 13).Value = 1314.6
$ws.Cells.Item(107, 14).Value = -4846.5

$ws.Cells.Item(134, 8).Value = 2678.0833
$ws.Cells.Item(134, 9).Value = 2200.7646
$ws.Cells.Item(134, 11).Value = 6602.293799999999
$ws.Cells.Item(134, 13).Value = -4067.293799999999

$ws.Cells.Item(136, 8).Value = 2612579.5
$ws.Cells.Item(136, 9).Value = 5216.3687
$ws.Cells.Item(136, 10).Value = 6423341
$ws.Cells.Item(136, 11).Value = 15649.1061
$ws.Cells.Item(136, 12).Value = 19270023
$ws.Cells.Item(136, 13).Value = -13099.1061
$ws.Cells.Item(136, 14).Value = -19275123

$ws.Cells.Item(141, 8).Value = 210652
$ws.Cells.Item(141, 10).Value = 224483.56
$ws.Cells.Item(141, 12).Value = 224483.56
$ws.Cells.Item(141, 14).Value = -234843.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 10737.637
$ws.Cells.Item(64, 10).Value = 12139.25
$ws.Cells.Item(64, 12).Value = 36417.75
$ws.Cells.Item(64, 14).Value = -36957.75

$ws.Cells.Item(67, 8).Value = 10737.637
$ws.Cells.Item(67, 10).Value = 12139.25
$ws.Cells.Item(67, 12).Value = 36417.75
$ws.Cells.Item(67, 14).Value = -38289.75

$ws.Cells.Item(86, 8).Value = 150.82353

$ws.Cells.Item(89, 8).Value = 150.82353

$ws.Cells.Item(107, 8).Value = 2110.96
$ws.Cells.Item(107, 10).Value = 2417.1904
$ws.Cells.Item(107, 12).Value = 7251.5712
$ws.Cells.Item(107, 14).Value = -11091.5712

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 4851.75
$ws.Cells.Item(21, 9).Value = 4858.091
$ws.Cells.Item(21, 11).Value = 4858.091
$ws.Cells.Item(21, 13).Value = -4685.091

$ws.Cells.Item(30, 8).Value = 4851.75
$ws.Cells.Item(30, 9).Value = 4858.091
$ws.Cells.Item(30, 11).Value = 4858.091
$ws.Cells.Item(30, 13).Value = -4753.091

$ws.Cells.Item(126, 8).Value = 12291.571
$ws.Cells.Item(126, 9).Value = 13756.833
$ws.Cells.Item(126, 11).Value = 41270.499
$ws.Cells.Item(126, 13).Value = -38800.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6613.3
$ws.Cells.Item(7, 9).Value = 7174
$ws.Cells.Item(7, 11).Value = 7174
$ws.Cells.Item(7, 13).Value = -7062

$ws.Cells.Item(23, 8).Value = 9103.214
$ws.Cells.Item(23, 10).Value = 8248.5
$ws.Cells.Item(23, 12).Value = 8248.5
$ws.Cells.Item(23, 14).Value = -8708.5

$ws.Cells.Item(40, 8).Value = 7965.2
$ws.Cells.Item(40, 9).Value = 7457.8335
$ws.Cells.Item(40, 11).Value = 7457.8335
$ws.Cells.Item(40, 13).Value = -7321.8335

$ws.Cells.Item(46, 8).Value = 4302.1113
$ws.Cells.Item(46, 9).Value = 1099.75
$ws.Cells.Item(46, 10).Value = 5217.0713
$ws.Cells.Item(46, 11).Value = 1099.75
$ws.Cells.Item(46, 12).Value = 5217.0713
$ws.Cells.Item(46, 13).Value = -911.75
$ws.Cells.Item(46, 14).Value = -5593.0713

$ws.Cells.Item(55, 8).Value = 1532.4073
$ws.Cells.Item(55, 9).Value = 1060.5294
$ws.Cells.Item(55, 10).Value = 2334.6
$ws.Cells.Item(55, 11).Value = 1060.5294
$ws.Cells.Item(55, 12).Value = 2334.6
$ws.Cells.Item(55, 13).Value = -887.5293999999999
$ws.Cells.Item(55, 14).Value = -2680.6

$ws.Cells.Item(100, 8).Value = 3499.7
$ws.Cells.Item(100, 9).Value = 2499.8333
$ws.Cells.Item(100, 10).Value = 4999.5
$ws.Cells.Item(100, 11).Value = 2499.8333
$ws.Cells.Item(100, 12).Value = 4999.5
$ws.Cells.Item(100, 13).Value = -1958.8333
$ws.Cells.Item(100, 14).Value = -6081.5

$ws.Cells.Item(122, 8).Value = 3052.7144
$ws.Cells.Item(122, 9).Value = 2847.3572
$ws.Cells.Item(122, 11).Value = 8542.071599999999
$ws.Cells.Item(122, 13).Value = -6092.071599999999

$ws.Cells.Item(126, 8).Value = 6613.3
$ws.Cells.Item(126, 9).Value = 7174
$ws.Cells.Item(126, 11).Value = 21522
$ws.Cells.Item(126, 13).Value = -19052

$ws.Cells.Item(136, 8).Value = 10424726
$ws.Cells.Item(136, 9).Value = 12503346
$ws.Cells.Item(136, 11).Value = 37510038
$ws.Cells.Item(136, 13).Value = -37507488

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2714.1428
$ws.Cells.Item(126, 9).Value = 2815.3845
$ws.Cells.Item(126, 11).Value = 8446.1535
$ws.Cells.Item(126, 13).Value = -5976.1535

$ws.Cells.Item(132, 8).Value = 4387434.5
$ws.Cells.Item(132, 9).Value = 4763209
$ws.Cells.Item(132, 10).Value = 3400
$ws.Cells.Item(132, 11).Value = 14289627
$ws.Cells.Item(132, 12).Value = 10200
$ws.Cells.Item(132, 13).Value = -14287097
$ws.Cells.Item(132, 14).Value = -15260
